$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 97, pushing the existing
# rows 97-172 down to 98-173 (dimension grows from T172 to T173).
$ws.Rows.Item(97).Insert()

$newRow = 97

$ws.Cells.Item($newRow, 1).Value2 = 9
$ws.Cells.Item($newRow, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value2 = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value2 = 44603
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value2 = 13
$ws.Cells.Item($newRow, 6).Value2 = "Fruta"
$ws.Cells.Item($newRow, 7).Value2 = 100101
$ws.Cells.Item($newRow, 8).Value2 = "Berries"
$ws.Cells.Item($newRow, 9).Value2 = 100101001
$ws.Cells.Item($newRow, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item($newRow, 11).Value2 = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value2 = "Primera"
$ws.Cells.Item($newRow, 13).Value2 = 350
$ws.Cells.Item($newRow, 14).Value2 = 3200
$ws.Cells.Item($newRow, 15).Value2 = 3200
$ws.Cells.Item($newRow, 16).Value2 = 3200
$ws.Cells.Item($newRow, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item($newRow, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value2 = 1600
$ws.Cells.Item($newRow, 20).Value2 = 2
